# Auto-generated Excel COM-interop script
# Applies value updates to Leve profit-tracking columns (H:N) across
# multiple crafting-job worksheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H2").Value2 = 301.1111
$ws.Range("I2").Value2 = 313.75
$ws.Range("K2").Value2 = 313.75
$ws.Range("M2").Value2 = -200.75

$ws.Range("H4").Value2 = 500
$ws.Range("I4").Value2 = 200
$ws.Range("J4").Value2 = 800
$ws.Range("K4").Value2 = 200
$ws.Range("L4").Value2 = 800
$ws.Range("M4").Value2 = -86
$ws.Range("N4").Value2 = -1028

$ws.Range("H5").Value2 = 243.16667
$ws.Range("I5").Value2 = 59.42857
$ws.Range("K5").Value2 = 59.42857
$ws.Range("M5").Value2 = 55.57143

$ws.Range("H9").Value2 = 190.16667
$ws.Range("I9").Value2 = 97.75
$ws.Range("J9").Value2 = 236.375
$ws.Range("K9").Value2 = 97.75
$ws.Range("L9").Value2 = 236.375
$ws.Range("M9").Value2 = 71.25
$ws.Range("N9").Value2 = -574.375

$ws.Range("H12").Value2 = 400
$ws.Range("I12").Value2 = 366.66666
$ws.Range("J12").Value2 = 500
$ws.Range("K12").Value2 = 366.66666
$ws.Range("L12").Value2 = 500
$ws.Range("M12").Value2 = -196.66666
$ws.Range("N12").Value2 = -840

$ws.Range("H19").Value2 = 744.6957
$ws.Range("I19").Value2 = 626.6667
$ws.Range("J19").Value2 = 966
$ws.Range("K19").Value2 = 626.6667
$ws.Range("L19").Value2 = 966
$ws.Range("M19").Value2 = -451.6667
$ws.Range("N19").Value2 = -1316

$ws.Range("H33").Value2 = 38.363636
$ws.Range("I33").Value2 = 39.2
$ws.Range("J33").Value2 = 30
$ws.Range("K33").Value2 = 39.2
$ws.Range("L33").Value2 = 30
$ws.Range("M33").Value2 = 189.8
$ws.Range("N33").Value2 = -488

$ws.Range("H40").Value2 = 2381.818
$ws.Range("I40").Value2 = 2193.6875
$ws.Range("J40").Value2 = 2883.5
$ws.Range("K40").Value2 = 2193.6875
$ws.Range("L40").Value2 = 2883.5
$ws.Range("M40").Value2 = -2018.6875
$ws.Range("N40").Value2 = -3233.5

$ws.Range("H76").Value2 = 4118211
$ws.Range("I76").Value2 = 4632504.5
$ws.Range("J76").Value2 = 3861.3333
$ws.Range("K76").Value2 = 4632504.5
$ws.Range("L76").Value2 = 3861.3333
$ws.Range("M76").Value2 = -4632189.5
$ws.Range("N76").Value2 = -4491.3333

$ws.Range("H79").Value2 = 4118211
$ws.Range("I79").Value2 = 4632504.5
$ws.Range("J79").Value2 = 3861.3333
$ws.Range("K79").Value2 = 4632504.5
$ws.Range("L79").Value2 = 3861.3333
$ws.Range("M79").Value2 = -4631412.5
$ws.Range("N79").Value2 = -6045.3333

$ws.Range("H107").Value2 = 654088.0600000001
$ws.Range("I107").Value2 = 794100.8
$ws.Range("J107").Value2 = 695.3333
$ws.Range("K107").Value2 = 794100.8
$ws.Range("L107").Value2 = 695.3333
$ws.Range("M107").Value2 = -792180.8
$ws.Range("N107").Value2 = -4535.3333

$ws.Range("H131").Value2 = 6691.25
$ws.Range("I131").Value2 = 1488.3334
$ws.Range("K131").Value2 = 4465.0002
$ws.Range("M131").Value2 = 574.9997999999996

$ws.Range("H138").Value2 = 4072223.8
$ws.Range("I138").Value2 = 1425061.9
$ws.Range("J138").Value2 = 5131088.5
$ws.Range("K138").Value2 = 4275185.699999999
$ws.Range("L138").Value2 = 15393265.5
$ws.Range("M138").Value2 = -4270045.699999999
$ws.Range("N138").Value2 = -15403545.5

$ws = $wb.Sheets("ARM")
$ws.Range("H4").Value2 = 287.16666
$ws.Range("I4").Value2 = 264.6
$ws.Range("J4").Value2 = 400
$ws.Range("K4").Value2 = 264.6
$ws.Range("L4").Value2 = 400
$ws.Range("M4").Value2 = -148.6
$ws.Range("N4").Value2 = -632

$ws.Range("H5").Value2 = 1000215.3
$ws.Range("I5").Value2 = 2000200.2
$ws.Range("J5").Value2 = 230.4
$ws.Range("K5").Value2 = 2000200.2
$ws.Range("L5").Value2 = 230.4
$ws.Range("M5").Value2 = -2000088.2
$ws.Range("N5").Value2 = -454.4

$ws.Range("H122").Value2 = 1677
$ws.Range("I122").Value2 = 1386.6364
$ws.Range("J122").Value2 = 2209.3333
$ws.Range("K122").Value2 = 4159.9092
$ws.Range("L122").Value2 = 6627.999899999999
$ws.Range("M122").Value2 = -1709.9092
$ws.Range("N122").Value2 = -11527.9999

$ws.Range("H132").Value2 = 2307.972
$ws.Range("I132").Value2 = 2068.3225
$ws.Range("J132").Value2 = 3958.889
$ws.Range("K132").Value2 = 6204.967500000001
$ws.Range("L132").Value2 = 11876.667
$ws.Range("M132").Value2 = -3674.967500000001
$ws.Range("N132").Value2 = -16936.667

$ws = $wb.Sheets("BSM")
$ws.Range("H4").Value2 = 1000215.3
$ws.Range("I4").Value2 = 2000200.2
$ws.Range("J4").Value2 = 230.4
$ws.Range("K4").Value2 = 2000200.2
$ws.Range("L4").Value2 = 230.4
$ws.Range("M4").Value2 = -2000085.2
$ws.Range("N4").Value2 = -460.4

$ws.Range("H107").Value2 = 898.8889
$ws.Range("J107").Value2 = 1098.2
$ws.Range("L107").Value2 = 1098.2
$ws.Range("N107").Value2 = -4938.2

$ws = $wb.Sheets("CRP")
$ws.Range("H4").Value2 = 1000000000
$ws.Range("J4").Value2 = 1000000000
$ws.Range("L4").Value2 = 1000000000
$ws.Range("N4").Value2 = -1000000224

$ws.Range("H7").Value2 = 120
$ws.Range("I7").Value2 = 140
$ws.Range("J7").Value2 = 104
$ws.Range("K7").Value2 = 140
$ws.Range("L7").Value2 = 104
$ws.Range("M7").Value2 = -27
$ws.Range("N7").Value2 = -330

$ws.Range("H134").Value2 = 2671
$ws.Range("I134").Value2 = 1444.9584
$ws.Range("K134").Value2 = 4334.8752
$ws.Range("M134").Value2 = -1799.8752

$ws.Range("H140").Value2 = 55000
$ws.Range("J140").Value2 = 55000
$ws.Range("L140").Value2 = 55000
$ws.Range("N140").Value2 = -65360

$ws = $wb.Sheets("CUL")
$ws.Range("H4").Value2 = 34169.1
$ws.Range("I4").Value2 = 137.4
$ws.Range("J4").Value2 = 68200.8
$ws.Range("K4").Value2 = 412.2
$ws.Range("L4").Value2 = 204602.4
$ws.Range("M4").Value2 = -300.2
$ws.Range("N4").Value2 = -204826.4

$ws.Range("H33").Value2 = 436.76923
$ws.Range("I33").Value2 = 230.33333
$ws.Range("J33").Value2 = 498.7
$ws.Range("K33").Value2 = 1381.99998
$ws.Range("L33").Value2 = 2992.2
$ws.Range("M33").Value2 = -1098.99998
$ws.Range("N33").Value2 = -3558.2

$ws.Range("H132").Value2 = 1229.6666
$ws.Range("I132").Value2 = 755.1667
$ws.Range("J132").Value2 = 1466.9166
$ws.Range("K132").Value2 = 6796.5003
$ws.Range("L132").Value2 = 13202.2494
$ws.Range("M132").Value2 = -4266.5003
$ws.Range("N132").Value2 = -18262.2494

$ws.Range("H137").Value2 = 7218708
$ws.Range("I137").Value2 = 14289555
$ws.Range("J137").Value2 = 147861.86
$ws.Range("K137").Value2 = 42868665
$ws.Range("L137").Value2 = 443585.58
$ws.Range("M137").Value2 = -42863565
$ws.Range("N137").Value2 = -453785.58

$ws = $wb.Sheets("GSM")
$ws.Range("H2").Value2 = 48.22222
$ws.Range("I2").Value2 = 44
$ws.Range("J2").Value2 = 50.333332
$ws.Range("K2").Value2 = 44
$ws.Range("L2").Value2 = 50.333332
$ws.Range("M2").Value2 = 69
$ws.Range("N2").Value2 = -276.333332

$ws.Range("H5").Value2 = 660
$ws.Range("I5").Value2 = 433.33334
$ws.Range("J5").Value2 = 1000
$ws.Range("K5").Value2 = 433.33334
$ws.Range("L5").Value2 = 1000
$ws.Range("M5").Value2 = -321.33334
$ws.Range("N5").Value2 = -1224

$ws.Range("H70").Value2 = 5384
$ws.Range("I70").Value2 = 5366.6943
$ws.Range("J70").Value2 = 5461.875
$ws.Range("K70").Value2 = 5366.6943
$ws.Range("L70").Value2 = 5461.875
$ws.Range("M70").Value2 = -5096.6943
$ws.Range("N70").Value2 = -6001.875

$ws.Range("H73").Value2 = 5384
$ws.Range("I73").Value2 = 5366.6943
$ws.Range("J73").Value2 = 5461.875
$ws.Range("K73").Value2 = 5366.6943
$ws.Range("L73").Value2 = 5461.875
$ws.Range("M73").Value2 = -4430.6943
$ws.Range("N73").Value2 = -7333.875

$ws.Range("H80").Value2 = 2388.125
$ws.Range("I80").Value2 = 2181
$ws.Range("J80").Value2 = 2733.3333
$ws.Range("K80").Value2 = 2181
$ws.Range("L80").Value2 = 2733.3333
$ws.Range("M80").Value2 = -1183
$ws.Range("N80").Value2 = -4729.3333

$ws.Range("H83").Value2 = 2388.125
$ws.Range("I83").Value2 = 2181
$ws.Range("J83").Value2 = 2733.3333
$ws.Range("K83").Value2 = 10905
$ws.Range("L83").Value2 = 13666.6665
$ws.Range("M83").Value2 = -5913
$ws.Range("N83").Value2 = -23650.6665

$ws.Range("H102").Value2 = 2263.25
$ws.Range("I102").Value2 = 1943.5834
$ws.Range("J102").Value2 = 3222.25
$ws.Range("K102").Value2 = 1943.5834
$ws.Range("L102").Value2 = 3222.25
$ws.Range("M102").Value2 = -321.5834
$ws.Range("N102").Value2 = -6466.25

$ws.Range("H122").Value2 = 1589230.1
$ws.Range("J122").Value2 = 2625
$ws.Range("L122").Value2 = 7875
$ws.Range("N122").Value2 = -12775

$ws.Range("H124").Value2 = 61000
$ws.Range("J124").Value2 = 61000
$ws.Range("L124").Value2 = 61000
$ws.Range("N124").Value2 = -70820

$ws.Range("H126").Value2 = 2447.9678
$ws.Range("I126").Value2 = 1818.8
$ws.Range("J126").Value2 = 2747.5715
$ws.Range("K126").Value2 = 5456.4
$ws.Range("L126").Value2 = 8242.7145
$ws.Range("M126").Value2 = -2986.4
$ws.Range("N126").Value2 = -13182.7145

$ws.Range("H128").Value2 = 54500
$ws.Range("J128").Value2 = 54500
$ws.Range("L128").Value2 = 54500
$ws.Range("N128").Value2 = -64460

$ws.Range("H132").Value2 = 4112.4546
$ws.Range("I132").Value2 = 3946.963
$ws.Range("J132").Value2 = 4857.1665
$ws.Range("K132").Value2 = 11840.889
$ws.Range("L132").Value2 = 14571.4995
$ws.Range("M132").Value2 = -9310.889000000001
$ws.Range("N132").Value2 = -19631.4995

$ws.Range("H133").Value2 = 18312.223
$ws.Range("J133").Value2 = 18312.223
$ws.Range("L133").Value2 = 18312.223
$ws.Range("N133").Value2 = -28432.223

$ws.Range("H135").Value2 = 1000000000
$ws.Range("J135").Value2 = 1000000000
$ws.Range("L135").Value2 = 1000000000
$ws.Range("N135").Value2 = -1000010140

$ws.Range("H140").Value2 = 0
$ws.Range("J140").Value2 = 0
$ws.Range("L140").Value2 = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Sheets("LTW")
$ws.Range("H7").Value2 = 2933.6667
$ws.Range("I7").Value2 = 1787.875
$ws.Range("J7").Value2 = 3506.5625
$ws.Range("K7").Value2 = 1787.875
$ws.Range("L7").Value2 = 3506.5625
$ws.Range("M7").Value2 = -1675.875
$ws.Range("N7").Value2 = -3730.5625

$ws.Range("H46").Value2 = 1045.6364
$ws.Range("I46").Value2 = 1085.7142
$ws.Range("J46").Value2 = 975.5
$ws.Range("K46").Value2 = 1085.7142
$ws.Range("L46").Value2 = 975.5
$ws.Range("M46").Value2 = -897.7141999999999
$ws.Range("N46").Value2 = -1351.5

$ws.Range("H126").Value2 = 2933.6667
$ws.Range("I126").Value2 = 1787.875
$ws.Range("J126").Value2 = 3506.5625
$ws.Range("K126").Value2 = 5363.625
$ws.Range("L126").Value2 = 10519.6875
$ws.Range("M126").Value2 = -2893.625
$ws.Range("N126").Value2 = -15459.6875

$ws = $wb.Sheets("WVR")
$ws.Range("H2").Value2 = 66686668
$ws.Range("I2").Value2 = 80010000
$ws.Range("K2").Value2 = 80010000
$ws.Range("M2").Value2 = -80009888

$ws.Range("H126").Value2 = 57776.555
$ws.Range("I126").Value2 = 85741.664
$ws.Range("J126").Value2 = 1846.3334
$ws.Range("K126").Value2 = 257224.992
$ws.Range("L126").Value2 = 5539.0002
$ws.Range("M126").Value2 = -254754.992
$ws.Range("N126").Value2 = -10479.0002
